$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 353 (pushes existing rows 353-480 down to 355-482),
# preserving formatting/styles of the rows below (as Excel's default Insert behavior does).
$ws.Rows.Item(353).Insert()
$ws.Rows.Item(353).Insert()

# New row 353: "Primera" quality entry for date 2022-09-05 (serial 44809)
$ws.Cells.Item(353,1).Value = 8
$ws.Cells.Item(353,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(353,3).Value = "Coquimbo"
$ws.Cells.Item(353,4).Value = 44809
$ws.Cells.Item(353,5).Value = 4
$ws.Cells.Item(353,6).Value = 100112017
$ws.Cells.Item(353,7).Value = "Apio"
$ws.Cells.Item(353,8).Value = "Americana (o)"
$ws.Cells.Item(353,9).Value = "Primera"
$ws.Cells.Item(353,10).Value = 2000
$ws.Cells.Item(353,11).Value = 9000
$ws.Cells.Item(353,12).Value = 10000
$ws.Cells.Item(353,13).Value = 9500
$ws.Cells.Item(353,14).Value = "`$/docena de matas"
$ws.Cells.Item(353,15).Value = "Provincia del Elquí"
$ws.Cells.Item(353,16).Value = 1583
$ws.Cells.Item(353,17).Value = 6
$ws.Cells.Item(353,18).Value = "Hortaliza"

# New row 354: "Segunda" quality entry for date 2022-09-05 (serial 44809)
$ws.Cells.Item(354,1).Value = 8
$ws.Cells.Item(354,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(354,3).Value = "Coquimbo"
$ws.Cells.Item(354,4).Value = 44809
$ws.Cells.Item(354,5).Value = 4
$ws.Cells.Item(354,6).Value = 100112017
$ws.Cells.Item(354,7).Value = "Apio"
$ws.Cells.Item(354,8).Value = "Americana (o)"
$ws.Cells.Item(354,9).Value = "Segunda"
$ws.Cells.Item(354,10).Value = 1400
$ws.Cells.Item(354,11).Value = 7000
$ws.Cells.Item(354,12).Value = 8000
$ws.Cells.Item(354,13).Value = 7500
$ws.Cells.Item(354,14).Value = "`$/docena de matas"
$ws.Cells.Item(354,15).Value = "Provincia del Elquí"
$ws.Cells.Item(354,16).Value = 1250
$ws.Cells.Item(354,17).Value = 6
$ws.Cells.Item(354,18).Value = "Hortaliza"
